# Actualizacion Datos Personales 4 nov
#
# The three summary sheets ("Estadisticos 1P", "Estadisticos 2P",
# "Estadisticos Final") each had an accidental duplicate of their second
# data row (row 3 == row 2). Remove the duplicate row in each.
#
# The "Rescatables" detail sheet had every student duplicated on two
# consecutive rows (except GARCIA). De-duplicate it down to one row per
# student, and add a new student (GONZALEZ MENDEZ CRISTIAN JAHIR) right
# after the first row.

$wb = $excel.ActiveWorkbook

# --- Estadisticos 1P / 2P / Final: drop the duplicated row 3 -----------
foreach ($name in @("Estadisticos 1P", "Estadisticos 2P", "Estadisticos Final")) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Rows.Item(3).Delete()
}

# --- Rescatables: de-duplicate rows, then insert the new student -------
$ws = $wb.Worksheets.Item("Rescatables")

$rowsToDelete = @(19, 17, 15, 13, 11, 9, 7, 3)
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}

$ws.Rows.Item(3).Insert()
$ws.Range("A3").Value = 20330051920081
$ws.Range("B3").Value = "GONZALEZ"
$ws.Range("C3").Value = "MENDEZ"
$ws.Range("D3").Value = "CRISTIAN JAHIR"
$ws.Range("E3").Value = "MANTIENE LOS MOTORES DE CA Y CC"
$ws.Range("F3").Value = "3AEV"
$ws.Range("G3").Value = 6
